# Apply the scraped-price refresh described in the commit:
# "Updated symbol list on Wed Jan 11 12:59:32 UTC 2023 with GitHub Actions"
#
# The source data keeps every cell as literal text (even price/volume
# numbers and percentages), so for any numeric-looking value we first pin
# the cell to a text NumberFormat ("@") before assigning .Value - otherwise
# Excel COM auto-coerces e.g. "277.69" to a Double or "1.12%" to 0.0112.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '277.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.12%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.24'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.25%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.13%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06406'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.69%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.009'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.23%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.202'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-6.90%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8869'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.05%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1526'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.54%'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07496'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.72%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02882'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.71%'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08976'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.91%'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001576'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.87%'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006398'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.90%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006109'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.64%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.476'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.57%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.303'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.18%'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.247'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.61%'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3085'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.03%'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.05181'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.30%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1350'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.37%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.924'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.18%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1517'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '9.91%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04398'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.39%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001175'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.37%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.003896'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.44%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-1.83%'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '1.53%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04100'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.86%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006803'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.30%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.46%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.78%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01165'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '4.38%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005332'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.09%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.628'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '9.30%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01851'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-7.46%'
